$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 'culture_collection' field/column (column T), per INSDC2017 review
$ws.Columns("T").Delete()
